$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (rows 133-134),
# pushing the existing rows 133:176 down to 135:178.
$ws.Range("A133:A134").EntireRow.Insert()

# New week's "Primera" quality row
$ws.Cells.Item(133,1).Value = 11
$ws.Cells.Item(133,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(133,3).Value = "Bíobío"
$ws.Cells.Item(133,4).Value = 44663
$ws.Cells.Item(133,5).Value = 8
$ws.Cells.Item(133,6).Value = 100112040
$ws.Cells.Item(133,7).Value = "Cilantro"
$ws.Cells.Item(133,8).Value = "Sin especificar"
$ws.Cells.Item(133,9).Value = "Primera"
$ws.Cells.Item(133,10).Value = 200
$ws.Cells.Item(133,11).Value = 600
$ws.Cells.Item(133,12).Value = 700
$ws.Cells.Item(133,13).Value = 650
$ws.Cells.Item(133,14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(133,15).Value = "Región de Ñuble"
$ws.Cells.Item(133,16).Value = 650
$ws.Cells.Item(133,17).Value = 1
$ws.Cells.Item(133,18).Value = "Hortaliza"

# New week's "Segunda" quality row
$ws.Cells.Item(134,1).Value = 11
$ws.Cells.Item(134,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(134,3).Value = "Bíobío"
$ws.Cells.Item(134,4).Value = 44663
$ws.Cells.Item(134,5).Value = 8
$ws.Cells.Item(134,6).Value = 100112040
$ws.Cells.Item(134,7).Value = "Cilantro"
$ws.Cells.Item(134,8).Value = "Sin especificar"
$ws.Cells.Item(134,9).Value = "Segunda"
$ws.Cells.Item(134,10).Value = 100
$ws.Cells.Item(134,11).Value = 500
$ws.Cells.Item(134,12).Value = 500
$ws.Cells.Item(134,13).Value = 500
$ws.Cells.Item(134,14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(134,15).Value = "Región de Ñuble"
$ws.Cells.Item(134,16).Value = 500
$ws.Cells.Item(134,17).Value = 1
$ws.Cells.Item(134,18).Value = "Hortaliza"
